$wb = $excel.ActiveWorkbook

# New values (per sheet) for rows 2..16 -> columns B (Cutoff) and C (Reaction_number)
$nbrB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$nbrC = @(704,709,705,699,694,692,685,664,661,654,654,650,647,617,616)

$barB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$barC = @(646,638,631,631,630,630,630,623,623,628,626,623,620,614,615)

$sheets = @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))
$dataB = @($nbrB, $barB)
$dataC = @($nbrC, $barC)

for ($s = 0; $s -lt $sheets.Count; $s++) {
    $ws = $sheets[$s]
    $bVals = $dataB[$s]
    $cVals = $dataC[$s]

    # Update existing rows 2..16 with new B/C values
    for ($i = 0; $i -lt $bVals.Count; $i++) {
        $row = $i + 2
        $ws.Range("B$row").Value = $bVals[$i]
        $ws.Range("C$row").Value = $cVals[$i]
    }

    # Remove the now-obsolete rows 17..20
    $ws.Rows("17:20").Delete()
}
